$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.355.00"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.425.37"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.07"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.48"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.516"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.64"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.13"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.03"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.803.83"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.405.36"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.209.45"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.25"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.35"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0920"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.78"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.26"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.53"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.58"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "49.02"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.83"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.21"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.125"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.21"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0767"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.43"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.86"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "125.80"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.54"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0290"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.928.60"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.11"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.64"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.82"

$ws.Range("E2").Value = "  +3.49%  "
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  +3.52%  "
$ws.Range("E6").Value = "  +5.29%  "
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +7.12%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("E18").Value = "  +3.13%  "
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E28").Value = "  -6.49%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("B32").Value = "Celestia"
$ws.Range("C32").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("E32").Value = "  +9.11%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E33").Value = "  +3.84%  "
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("E43").Value = "  -3.39%  "
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("E48").Value = "  +15.34%  "
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("E50").Value = "  +5.45%  "
$ws.Range("E51").Value = "  +1.97%  "
